$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "01/10/2018"
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 4).Value = 50.82170884740103

$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = "01/01/2019"
$ws.Cells.Item(3, 3).Style = "Normal"
$ws.Cells.Item(3, 4).Value = 50.87924241045994

$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "01/04/2019"
$ws.Cells.Item(4, 3).Style = "Normal"
$ws.Cells.Item(4, 4).Value = 51.20988551004415

$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "01/07/2019"
$ws.Cells.Item(5, 3).Style = "Normal"
$ws.Cells.Item(5, 4).Value = 51.28578105476016

$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "01/10/2019"
$ws.Cells.Item(6, 3).Style = "Normal"
$ws.Cells.Item(6, 4).Value = 51.13267992212379

$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "01/01/2020"
$ws.Cells.Item(7, 3).Style = "Normal"
$ws.Cells.Item(7, 4).Value = 50.48747113657745

$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "01/04/2020"
$ws.Cells.Item(8, 3).Style = "Normal"
$ws.Cells.Item(8, 4).Value = 46.13243293229446

$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = "01/07/2020"
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).Value = 46.40541129024624

$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "01/10/2020"
$ws.Cells.Item(10, 3).Style = "Normal"
$ws.Cells.Item(10, 4).Value = 48.02080774100882

$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "01/01/2021"
$ws.Cells.Item(11, 3).Style = "Normal"
$ws.Cells.Item(11, 4).Value = 48.26400679117148

$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = "01/04/2021"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).Value = 49.05990792087598

$ws.Cells.Item(13, 3).NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = "01/07/2021"
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).Value = 50.01221758580504

$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "01/10/2021"
$ws.Cells.Item(14, 3).Style = "Normal"
$ws.Cells.Item(14, 4).Value = 50.54599696042929

$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "01/01/2022"
$ws.Cells.Item(15, 3).Style = "Normal"
$ws.Cells.Item(15, 4).Value = 50.20649354297969

$ws.Cells.Item(16, 3).NumberFormat = "@"
$ws.Cells.Item(16, 3).Value = "01/04/2022"
$ws.Cells.Item(16, 3).Style = "Normal"
$ws.Cells.Item(16, 4).Value = 50.644573244835

$ws.Cells.Item(17, 3).NumberFormat = "@"
$ws.Cells.Item(17, 3).Value = "01/07/2022"
$ws.Cells.Item(17, 3).Style = "Normal"
$ws.Cells.Item(17, 4).Value = 50.73397663219978

$ws.Cells.Item(18, 3).NumberFormat = "@"
$ws.Cells.Item(18, 3).Value = "01/10/2022"
$ws.Cells.Item(18, 3).Style = "Normal"
$ws.Cells.Item(18, 4).Value = 50.28041736538103

$ws.Cells.Item(19, 3).NumberFormat = "@"
$ws.Cells.Item(19, 3).Value = "01/01/2023"
$ws.Cells.Item(19, 3).Style = "Normal"
$ws.Cells.Item(19, 4).Value = 49.8763055346298

$ws.Cells.Item(20, 3).NumberFormat = "@"
$ws.Cells.Item(20, 3).Value = "01/04/2023"
$ws.Cells.Item(20, 3).Style = "Normal"
$ws.Cells.Item(20, 4).Value = 49.93175740919557

$ws.Cells.Item(21, 3).NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = "01/07/2023"
$ws.Cells.Item(21, 3).Style = "Normal"
$ws.Cells.Item(21, 4).Value = 50.12536671502129

$ws.Cells.Item(22, 3).NumberFormat = "@"
$ws.Cells.Item(22, 3).Value = "01/10/2023"
$ws.Cells.Item(22, 3).Style = "Normal"
$ws.Cells.Item(22, 4).Value = 50.46478162897979

$ws.Cells.Item(23, 3).NumberFormat = "@"
$ws.Cells.Item(23, 3).Value = "01/10/2018"
$ws.Cells.Item(23, 3).Style = "Normal"
$ws.Cells.Item(23, 4).Value = 44.14071449972682

$ws.Cells.Item(24, 3).NumberFormat = "@"
$ws.Cells.Item(24, 3).Value = "01/01/2019"
$ws.Cells.Item(24, 3).Style = "Normal"
$ws.Cells.Item(24, 4).Value = 43.96550206811581

$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 3).Value = "01/04/2019"
$ws.Cells.Item(25, 3).Style = "Normal"
$ws.Cells.Item(25, 4).Value = 44.24777205533388

$ws.Cells.Item(26, 3).NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = "01/07/2019"
$ws.Cells.Item(26, 3).Style = "Normal"
$ws.Cells.Item(26, 4).Value = 44.27124475573578

$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = "01/10/2019"
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(27, 4).Value = 44.20076434907612

$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "01/01/2020"
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(28, 4).Value = 43.72505077386371

$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = "01/04/2020"
$ws.Cells.Item(29, 3).Style = "Normal"
$ws.Cells.Item(29, 4).ClearContents()

$ws.Cells.Item(30, 3).NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = "01/07/2020"
$ws.Cells.Item(30, 3).Style = "Normal"

$ws.Cells.Item(31, 3).NumberFormat = "@"
$ws.Cells.Item(31, 3).Value = "01/10/2020"
$ws.Cells.Item(31, 3).Style = "Normal"

$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = "01/01/2021"
$ws.Cells.Item(32, 3).Style = "Normal"

$ws.Cells.Item(33, 3).NumberFormat = "@"
$ws.Cells.Item(33, 3).Value = "01/04/2021"
$ws.Cells.Item(33, 3).Style = "Normal"

$ws.Cells.Item(34, 3).NumberFormat = "@"
$ws.Cells.Item(34, 3).Value = "01/07/2021"
$ws.Cells.Item(34, 3).Style = "Normal"

$ws.Cells.Item(35, 3).NumberFormat = "@"
$ws.Cells.Item(35, 3).Value = "01/10/2021"
$ws.Cells.Item(35, 3).Style = "Normal"

$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = "01/01/2022"
$ws.Cells.Item(36, 3).Style = "Normal"

$ws.Cells.Item(37, 3).NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = "01/04/2022"
$ws.Cells.Item(37, 3).Style = "Normal"
$ws.Cells.Item(37, 4).Value = 43.73755734441271

$ws.Cells.Item(38, 3).NumberFormat = "@"
$ws.Cells.Item(38, 3).Value = "01/07/2022"
$ws.Cells.Item(38, 3).Style = "Normal"
$ws.Cells.Item(38, 4).Value = 43.96375834284331

$ws.Cells.Item(39, 3).NumberFormat = "@"
$ws.Cells.Item(39, 3).Value = "01/10/2022"
$ws.Cells.Item(39, 3).Style = "Normal"
$ws.Cells.Item(39, 4).Value = 43.55560160958845

$ws.Cells.Item(40, 3).NumberFormat = "@"
$ws.Cells.Item(40, 3).Value = "01/01/2023"
$ws.Cells.Item(40, 3).Style = "Normal"
$ws.Cells.Item(40, 4).Value = 42.77828569457142

$ws.Cells.Item(41, 3).NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = "01/04/2023"
$ws.Cells.Item(41, 3).Style = "Normal"
$ws.Cells.Item(41, 4).Value = 42.80053067764779

$ws.Cells.Item(42, 3).NumberFormat = "@"
$ws.Cells.Item(42, 3).Value = "01/07/2023"
$ws.Cells.Item(42, 3).Style = "Normal"
$ws.Cells.Item(42, 4).Value = 43.46097581963686

$ws.Cells.Item(43, 3).NumberFormat = "@"
$ws.Cells.Item(43, 3).Value = "01/10/2023"
$ws.Cells.Item(43, 3).Style = "Normal"
$ws.Cells.Item(43, 4).Value = 43.34118375135377

$ws.Cells.Item(44, 3).NumberFormat = "@"
$ws.Cells.Item(44, 3).Value = "01/10/2018"
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).Value = 46.58792650918635

$ws.Cells.Item(45, 3).NumberFormat = "@"
$ws.Cells.Item(45, 3).Value = "01/01/2019"
$ws.Cells.Item(45, 3).Style = "Normal"
$ws.Cells.Item(45, 4).Value = 46.39895242252292

$ws.Cells.Item(46, 3).NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = "01/04/2019"
$ws.Cells.Item(46, 3).Style = "Normal"
$ws.Cells.Item(46, 4).Value = 47.95296167247386

$ws.Cells.Item(47, 3).NumberFormat = "@"
$ws.Cells.Item(47, 3).Value = "01/07/2019"
$ws.Cells.Item(47, 3).Style = "Normal"
$ws.Cells.Item(47, 4).Value = 47.63146458061712

$ws.Cells.Item(48, 3).NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = "01/10/2019"
$ws.Cells.Item(48, 3).Style = "Normal"
$ws.Cells.Item(48, 4).Value = 48.78577623590633

$ws.Cells.Item(49, 3).NumberFormat = "@"
$ws.Cells.Item(49, 3).Value = "01/01/2020"
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = 46.8628299437473

$ws.Cells.Item(50, 3).NumberFormat = "@"
$ws.Cells.Item(50, 3).Value = "01/04/2020"
$ws.Cells.Item(50, 3).Style = "Normal"
$ws.Cells.Item(50, 4).ClearContents()

$ws.Cells.Item(51, 3).NumberFormat = "@"
$ws.Cells.Item(51, 3).Value = "01/07/2020"
$ws.Cells.Item(51, 3).Style = "Normal"

$ws.Cells.Item(52, 3).NumberFormat = "@"
$ws.Cells.Item(52, 3).Value = "01/10/2020"
$ws.Cells.Item(52, 3).Style = "Normal"

$ws.Cells.Item(53, 3).NumberFormat = "@"
$ws.Cells.Item(53, 3).Value = "01/01/2021"
$ws.Cells.Item(53, 3).Style = "Normal"

$ws.Cells.Item(54, 3).NumberFormat = "@"
$ws.Cells.Item(54, 3).Value = "01/04/2021"
$ws.Cells.Item(54, 3).Style = "Normal"

$ws.Cells.Item(55, 3).NumberFormat = "@"
$ws.Cells.Item(55, 3).Value = "01/07/2021"
$ws.Cells.Item(55, 3).Style = "Normal"

$ws.Cells.Item(56, 3).NumberFormat = "@"
$ws.Cells.Item(56, 3).Value = "01/10/2021"
$ws.Cells.Item(56, 3).Style = "Normal"

$ws.Cells.Item(57, 3).NumberFormat = "@"
$ws.Cells.Item(57, 3).Value = "01/01/2022"
$ws.Cells.Item(57, 3).Style = "Normal"

$ws.Cells.Item(58, 3).NumberFormat = "@"
$ws.Cells.Item(58, 3).Value = "01/04/2022"
$ws.Cells.Item(58, 3).Style = "Normal"
$ws.Cells.Item(58, 4).Value = 46.62420382165605

$ws.Cells.Item(59, 3).NumberFormat = "@"
$ws.Cells.Item(59, 3).Value = "01/07/2022"
$ws.Cells.Item(59, 3).Style = "Normal"
$ws.Cells.Item(59, 4).Value = 45.76271186440678

$ws.Cells.Item(60, 3).NumberFormat = "@"
$ws.Cells.Item(60, 3).Value = "01/10/2022"
$ws.Cells.Item(60, 3).Style = "Normal"
$ws.Cells.Item(60, 4).Value = 46.36209813874789

$ws.Cells.Item(61, 3).NumberFormat = "@"
$ws.Cells.Item(61, 3).Value = "01/01/2023"
$ws.Cells.Item(61, 3).Style = "Normal"
$ws.Cells.Item(61, 4).Value = 44.61798227100042

$ws.Cells.Item(62, 3).NumberFormat = "@"
$ws.Cells.Item(62, 3).Value = "01/04/2023"
$ws.Cells.Item(62, 3).Style = "Normal"
$ws.Cells.Item(62, 4).Value = 44.43976411120472

$ws.Cells.Item(63, 3).NumberFormat = "@"
$ws.Cells.Item(63, 3).Value = "01/07/2023"
$ws.Cells.Item(63, 3).Style = "Normal"
$ws.Cells.Item(63, 4).Value = 45.50042052144659

$ws.Cells.Item(64, 3).NumberFormat = "@"
$ws.Cells.Item(64, 3).Value = "01/10/2023"
$ws.Cells.Item(64, 3).Style = "Normal"
$ws.Cells.Item(64, 4).Value = 44.94334872010072
